$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Names
$ws.Range("C2").Value = "Hartmut"

# B3 holds a card-number-like string that must stay text (same style as before).
# Force text storage via the "@" text format, then restore the original cell
# style (s=8) by copying formatting back from a neighboring cell that already
# uses that style, so only the value/content changes.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("C3").Value = "Mohaupt"

# Statement start balance line
$ws.Range("D5").Value = "KONTOSTAND AM 16.07.2025"

# Row 6
$ws.Range("B6").Value = "17.07."
$ws.Range("C6").Value = "18.07."
$ws.Range("D6").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E6").Value = "24,58-"

# Row 7
$ws.Range("B7").Value = "20.07."
$ws.Range("C7").Value = "21.07."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-70414131"
$ws.Range("E7").Value = "53,65-"

# Row 8
$ws.Range("B8").Value = "22.07."
$ws.Range("C8").Value = "23.07."
$ws.Range("D8").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E8").Value = "69,28-"

# Statement end balance line
$ws.Range("D12").Value = "KONTOSTAND AM 25.07.2025"
$ws.Range("E12").Value = "147,51-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 30.07.2025"
